$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values so they are not
# auto-converted to real numbers (source data stores these as text,
# same as every other cell in the Price column).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '29.308.72'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.873.49'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '0.7084'
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("D6").Value = '241.92'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.07801'
$ws.Range("E8").Value = '  +1.08%  '
$ws.Range("D9").Value = '0.3108'
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").Value = '25.09'
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").Value = '0.08376'
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '1.877.00'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '5.235'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '0.7175'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").Value = '0.000008391'
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("D17").Value = '6.149'
$ws.Range("E17").Value = '  +3.17%  '
$ws.Range("D18").Value = '29.313.59'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '240.67'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '13.20'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.122.24'
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").Value = '0.9999'
$ws.Range("E23").Value = '  -1.91%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '0.1599'
$ws.Range("E25").Value = '  -1.96%  '
$ws.Range("D26").Value = '162.79'
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("D27").Value = '9.039'
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").Value = '18.50'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").Value = '1.504'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '4.414'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = '4.340'
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("D32").Value = '1.222'
$ws.Range("E32").Value = '  -4.82%  '
$ws.Range("D33").Value = '0.05352'
$ws.Range("E33").Value = '  +2.06%  '
$ws.Range("D34").Value = '1.942'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").Value = '1.174'
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").Value = '0.7461'
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("D37").Value = '2.685'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("D39").Value = '1.241.34'
$ws.Range("E39").Value = '  +7.01%  '
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '6.507'
$ws.Range("E41").Value = '  +2.28%  '
$ws.Range("D42").Value = '0.8932'
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("D43").Value = '109.91'
$ws.Range("E43").Value = '  +4.86%  '
$ws.Range("D44").Value = '72.29'
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D45").Value = '0.00000000131'
$ws.Range("E45").Value = '  +12.30%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("D48").Value = '0.5193'
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("D49").Value = '1.794'
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("D50").Value = '9.458'
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("D51").Value = '0.4338'
$ws.Range("E51").Value = '  +0.65%  '
